$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Flags")
$ws2 = $wb.Worksheets.Item("Tests")

# --- "Tests" sheet: retire the stale last test row (row 42) ---
# Re-stamp A42/D42 with the "normal" row formatting (copy format only,
# content/formula untouched) before clearing out the now-obsolete
# sample data in B42:D42.
$ws2.Range("A2").Copy()
$ws2.Range("A42").PasteSpecial(-4122)
$ws2.Range("D2").Copy()
$ws2.Range("D42").PasteSpecial(-4122)

$ws2.Range("B42").ClearContents()
$ws2.Range("C42").ClearContents()
$ws2.Range("D42").ClearContents()

# --- "Flags" sheet: Categories -> Debug, AllColors -> False ---
$ws1.Range("B3").Value = "Debug"

# Assign "False" as literal text (not the Boolean FALSE) by writing it as
# a formula result and then freezing that result into a plain value with
# a values-only paste.
$ws1.Range("B4").Formula = "=""False"""
$ws1.Range("B4").Copy()
$ws1.Range("B4").PasteSpecial(-4163)

# Drop the explicit row height override on row 4 now that it no longer
# holds a tall multi-line description.
$ws1.Rows.Item(4).AutoFit()

# --- Restore cursor/selection state ---
[void]$ws1.Range("B5").Select()
$ws2.Activate()
[void]$ws2.Range("J18").Select()

$wb.Saved = $true
